$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R:T").Delete() | Out-Null
$ws.Range("O:P").Delete() | Out-Null
